# Apply the edits described by the commit:
# "Created Function for Gaussian Quadrature Scheme, and exported it to the
#  Averaged Intensities files."
#
# 1. Rename the worksheet (and thus the workbook's sheet entry) from
#    "BrassA-HW30.xpc" to "BrassA".
# 2. Tweak a few values in row 13 / row 15 that changed by floating point
#    noise after recomputation with the new Gaussian quadrature routine.
# 3. Append a new row (row 16) with a fresh set of averaged-intensity
#    results (HKL index 14, label "HexGrid-60degTilt5degRes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename sheet -------------------------------------------------
$ws.Name = "BrassA"

# --- 2. Small floating point corrections in existing rows -------------
$ws.Range("F13").Value = 0.9943235584194052
$ws.Range("L13").Value = 0.9940080201044336
$ws.Range("O13").Value = 0.9946842827013727
$ws.Range("J15").Value = 0.998407628726652

# --- 3. Append new row 16 ---------------------------------------------
$ws.Cells.Item(16, 1).Value = 14
# Carry over the bold/bordered/centered style used by the rest of column A
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 1.034373068099897
$ws.Cells.Item(16, 4).Value = 0.9546904888427717
$ws.Cells.Item(16, 5).Value = 0.9948365564448869
$ws.Cells.Item(16, 6).Value = 0.9842844965884121
$ws.Cells.Item(16, 7).Value = 1.034373068099897
$ws.Cells.Item(16, 8).Value = 0.9546904888427717
$ws.Cells.Item(16, 9).Value = 1.009099559340002
$ws.Cells.Item(16, 10).Value = 0.9803122053658803
$ws.Cells.Item(16, 11).Value = 1.008112316780786
$ws.Cells.Item(16, 12).Value = 0.9676901038902267
$ws.Cells.Item(16, 13).Value = 1.034373068099897
$ws.Cells.Item(16, 14).Value = 0.9747635226438294
$ws.Cells.Item(16, 15).Value = 0.9920461524939921
$ws.Cells.Item(16, 16).Value = 0.9916748494191079
